$d = $word.ActiveDocument

# The document's first paragraph holds the "**ID__...__ID**" merge-field
# placeholder. The edit:
#   1. Gives the paragraph a border whose only effect is 5-twip padding
#      on all four sides (no visible rule -- i.e. <w:pBdr><w:top w:space="5"/>...).
#   2. Increases the paragraph's left indent from 120 to 225 twips.
#   3. Renames the placeholder id text.
#   4. Removes the trailing single-space run that used to follow the id.

$p1 = $d.Paragraphs(1)

# --- 1) Paragraph border / padding -----------------------------------
$borders = $p1.Range.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- 2) Left indent ----------------------------------------------------
# 225 twips == 11.25 points (Word's ParagraphFormat.LeftIndent is in points).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- 3 & 4) Replace the id text and drop the trailing space run -------
$oldId = "**ID__AFFARS_pgi_5335_topic_2__ID**"
$newId = "**ID__AFFARS_AFICC_PGI_5335__ID**"

$paraRange = $p1.Range
$idStart = $paraRange.Start
$idEnd = $idStart + $oldId.Length

# Delete everything between the end of the id text and the paragraph
# mark (i.e. the lone trailing space / its run) first, so the id range's
# offsets remain valid afterwards.
$trailRange = $d.Range($idEnd, $paraRange.End - 1)
if ($trailRange.Start -lt $trailRange.End) {
    $trailRange.Delete()
}

$idRange = $d.Range($idStart, $idEnd)
$idRange.Text = $newId
